$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.963.13"
$ws.Cells.Item(2, 5).Value = "  +1.48%  "
$ws.Cells.Item(3, 4).Value = "3.902.89"
$ws.Cells.Item(3, 5).Value = "  -0.23%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "483.48"
$ws.Cells.Item(5, 5).Value = "  +2.71%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "145.62"
$ws.Cells.Item(6, 5).Value = "  +0.67%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.624"
$ws.Cells.Item(7, 5).Value = "  -1.00%  "
$ws.Cells.Item(8, 5).Value = "  -0.03%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.726"
$ws.Cells.Item(9, 5).Value = "  -2.41%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.168"
$ws.Cells.Item(10, 5).Value = "  +2.23%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0000361"
$ws.Cells.Item(11, 5).Value = "  +8.93%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "42.60"
$ws.Cells.Item(12, 5).Value = "  -1.78%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.70"
$ws.Cells.Item(13, 5).Value = "  +2.76%  "
$ws.Cells.Item(14, 4).Value = "4.525.45"
$ws.Cells.Item(14, 5).Value = "  -0.31%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "14.83"
$ws.Cells.Item(15, 5).Value = "  -2.46%  "
$ws.Cells.Item(16, 4).Value = "3.917.13"
$ws.Cells.Item(16, 5).Value = "  +0.61%  "
$ws.Cells.Item(17, 5).Value = "  -0.20%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "19.74"
$ws.Cells.Item(18, 5).Value = "  -1.96%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.13"
$ws.Cells.Item(19, 5).Value = "  -3.09%  "
$ws.Cells.Item(20, 4).Value = "68.058.53"
$ws.Cells.Item(20, 5).Value = "  +1.17%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "448.48"
$ws.Cells.Item(21, 5).Value = "  +3.94%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "14.83"
$ws.Cells.Item(22, 5).Value = "  -0.39%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.36"
$ws.Cells.Item(23, 5).Value = "  +0.70%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "88.96"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "11.66"
$ws.Cells.Item(25, 5).Value = "  +15.14%  "
$ws.Cells.Item(26, 2).Value = "PancakeSwap"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.57"
$ws.Cells.Item(26, 5).Value = "  -0.62%  "
$ws.Cells.Item(27, 2).Value = "RenderToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.59"
$ws.Cells.Item(27, 5).Value = "  +8.30%  "
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "38.54"
$ws.Cells.Item(28, 5).Value = "  -0.75%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.82"
$ws.Cells.Item(29, 5).Value = "  +3.11%  "
$ws.Cells.Item(30, 2).Value = "Bittensor"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "690.05"
$ws.Cells.Item(30, 5).Value = "  -5.89%  "
$ws.Cells.Item(31, 2).Value = "Cosmos"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "13.32"
$ws.Cells.Item(31, 5).Value = "  -3.01%  "
$ws.Cells.Item(32, 5).Value = "  -1.82%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.85"
$ws.Cells.Item(33, 5).Value = "  +2.39%  "
$ws.Cells.Item(34, 4).Value = "0.0₃0961"
$ws.Cells.Item(34, 5).Value = "  +30.34%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "41.39"
$ws.Cells.Item(35, 5).Value = "  -6.46%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "58.66"
$ws.Cells.Item(36, 5).Value = "  +0.61%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.71"
$ws.Cells.Item(37, 5).Value = "  +6.25%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.150"
$ws.Cells.Item(38, 5).Value = "  -6.11%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.999"
$ws.Cells.Item(39, 5).Value = "  +0.00%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0474"
$ws.Cells.Item(40, 5).Value = "  -1.04%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.76"
$ws.Cells.Item(41, 5).Value = "  +12.01%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.03"
$ws.Cells.Item(42, 5).Value = "  +7.84%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.02"
$ws.Cells.Item(43, 5).Value = "  -5.53%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.349"
$ws.Cells.Item(44, 5).Value = "  +2.92%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.142"
$ws.Cells.Item(45, 5).Value = "  +0.17%  "
$ws.Cells.Item(46, 5).Value = "  -0.11%  "
$ws.Cells.Item(47, 5).Value = "  -1.31%  "
$ws.Cells.Item(48, 2).Value = "ARBITRUM"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.13"
$ws.Cells.Item(48, 5).Value = "  -3.29%  "
$ws.Cells.Item(49, 2).Value = "Monero"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "145.79"
$ws.Cells.Item(49, 5).Value = "  +1.58%  "
$ws.Cells.Item(50, 2).Value = "FLOKI"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.000269"
$ws.Cells.Item(50, 5).Value = "  +69.28%  "
$ws.Cells.Item(51, 2).Value = "ApeXProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "3.10"
$ws.Cells.Item(51, 5).Value = "  -3.65%  "
